$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.296.44"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.864.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.17%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'234.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.76%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  +0.08%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4691"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.19%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -1.99%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06559"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.24%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'21.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -2.28%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07822"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.33%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'WrappedEther"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'1.886.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.12%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'Litecoin"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'96.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.08%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.6962"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.32%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'5.091"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.11%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'267.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.59%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'30.389.93"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.50%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'13.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.04%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.000007613"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +2.41%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +0.07%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'2.127.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.67%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.05%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.214"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.37%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -0.43%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'9.436"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.48%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'166.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.20%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'18.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.20%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'1.937"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.71%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -1.79%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.09907"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.61%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'4.349"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.38%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'1.455"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -1.09%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +0.02%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.04722"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'1.131"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.19%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.7016"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.19%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +0.52%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.01871"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.34%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  +5.32%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.45%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'72.52"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.94%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'1.947"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +0.03%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.4170"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +0.35%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +0.15%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.8350"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'103.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.24%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'969.63"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +2.07%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'7.099"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.90%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'9.143"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.19%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'34.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.84%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.05681"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.39%  "
$ws.Range("E51").Style = "Normal"
